$wb = $excel.ActiveWorkbook

# Sheet ALC row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 374.84616
$ws.Cells.Item(2, 10).Value = 549.75
$ws.Cells.Item(2, 12).Value = 549.75
$ws.Cells.Item(2, 14).Value = -775.75

# Sheet ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3082.9348
$ws.Cells.Item(17, 10).Value = 3034.7778
$ws.Cells.Item(17, 12).Value = 9104.3334
$ws.Cells.Item(17, 14).Value = -9440.3334

# Sheet ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1073.1765
$ws.Cells.Item(19, 9).Value = 1196.3334
$ws.Cells.Item(19, 10).Value = 934.625
$ws.Cells.Item(19, 11).Value = 1196.3334
$ws.Cells.Item(19, 12).Value = 934.625
$ws.Cells.Item(19, 13).Value = -1021.3334
$ws.Cells.Item(19, 14).Value = -1284.625

# Sheet ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3795.92
$ws.Cells.Item(40, 10).Value = 3909
$ws.Cells.Item(40, 12).Value = 3909
$ws.Cells.Item(40, 14).Value = -4259

# Sheet ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1781.0834
$ws.Cells.Item(98, 9).Value = 1488.7273
$ws.Cells.Item(98, 11).Value = 1488.7273
$ws.Cells.Item(98, 13).Value = 9.272699999999986

# Sheet ALC row 104
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(104, 8).Value = 620.4
$ws.Cells.Item(104, 10).Value = 1937
$ws.Cells.Item(104, 12).Value = 5811
$ws.Cells.Item(104, 14).Value = -9305

# Sheet ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 879854.2
$ws.Cells.Item(112, 10).Value = 1152117.8
$ws.Cells.Item(112, 12).Value = 3456353.4
$ws.Cells.Item(112, 14).Value = -3458569.4

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 12845.167
$ws.Cells.Item(116, 9).Value = 14840.941
$ws.Cells.Item(116, 10).Value = 7998.2856
$ws.Cells.Item(116, 11).Value = 14840.941
$ws.Cells.Item(116, 12).Value = 7998.2856
$ws.Cells.Item(116, 13).Value = -11398.941
$ws.Cells.Item(116, 14).Value = -14882.2856

# Sheet ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(117, 8).Value = 86139.375
$ws.Cells.Item(117, 10).Value = 86139.375
$ws.Cells.Item(117, 12).Value = 86139.375
$ws.Cells.Item(117, 14).Value = -95317.375

# Sheet ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1781.0834
$ws.Cells.Item(122, 9).Value = 1488.7273
$ws.Cells.Item(122, 11).Value = 4466.1819
$ws.Cells.Item(122, 13).Value = -2016.1819

# Sheet ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 4112.5
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 4112.5
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 13).ClearContents()
$ws.Cells.Item(125, 14).Value = -41932.5

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1647.7727
$ws.Cells.Item(132, 9).Value = 1588
$ws.Cells.Item(132, 11).Value = 4764
$ws.Cells.Item(132, 13).Value = -2234

# Sheet ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 3249.0833
$ws.Cells.Item(135, 9).Value = 1398.7142
$ws.Cells.Item(135, 11).Value = 12588.4278
$ws.Cells.Item(135, 13).Value = -10053.4278

# Sheet ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 718.2857
$ws.Cells.Item(5, 9).Value = 718.2857
$ws.Cells.Item(5, 11).Value = 718.2857
$ws.Cells.Item(5, 13).Value = -606.2857

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3785.2856
$ws.Cells.Item(45, 10).Value = 4249.5
$ws.Cells.Item(45, 12).Value = 4249.5
$ws.Cells.Item(45, 14).Value = -5003.5

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2303.8125
$ws.Cells.Item(74, 9).Value = 2255.9312
$ws.Cells.Item(74, 11).Value = 2255.9312
$ws.Cells.Item(74, 13).Value = -1381.9312

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2303.8125
$ws.Cells.Item(77, 9).Value = 2255.9312
$ws.Cells.Item(77, 11).Value = 11279.656
$ws.Cells.Item(77, 13).Value = -6911.655999999999

# Sheet ARM row 81
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(81, 8).Value = 125997.5
$ws.Cells.Item(81, 10).Value = 125997.5
$ws.Cells.Item(81, 12).Value = 125997.5
$ws.Cells.Item(81, 14).Value = -127993.5

# Sheet ARM row 84
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(84, 8).Value = 125997.5
$ws.Cells.Item(84, 10).Value = 125997.5
$ws.Cells.Item(84, 12).Value = 377992.5
$ws.Cells.Item(84, 14).Value = -387976.5

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3286.6155
$ws.Cells.Item(122, 10).Value = 3485.25
$ws.Cells.Item(122, 12).Value = 10455.75
$ws.Cells.Item(122, 14).Value = -15355.75

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2481.4238
$ws.Cells.Item(132, 9).Value = 2501.94
$ws.Cells.Item(132, 10).Value = 2367.4443
$ws.Cells.Item(132, 11).Value = 7505.82
$ws.Cells.Item(132, 12).Value = 7102.3329
$ws.Cells.Item(132, 13).Value = -4975.82
$ws.Cells.Item(132, 14).Value = -12162.3329

# Sheet BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 718.2857
$ws.Cells.Item(4, 9).Value = 718.2857
$ws.Cells.Item(4, 11).Value = 718.2857
$ws.Cells.Item(4, 13).Value = -603.2857

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2001.1875
$ws.Cells.Item(86, 9).Value = 2522.5
$ws.Cells.Item(86, 10).Value = 1132.3334
$ws.Cells.Item(86, 11).Value = 2522.5
$ws.Cells.Item(86, 12).Value = 1132.3334
$ws.Cells.Item(86, 13).Value = -1399.5
$ws.Cells.Item(86, 14).Value = -3378.3334

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2001.1875
$ws.Cells.Item(89, 9).Value = 2522.5
$ws.Cells.Item(89, 10).Value = 1132.3334
$ws.Cells.Item(89, 11).Value = 12612.5
$ws.Cells.Item(89, 12).Value = 5661.666999999999
$ws.Cells.Item(89, 13).Value = -6996.5
$ws.Cells.Item(89, 14).Value = -16893.667

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3921
$ws.Cells.Item(31, 9).Value = 2436.2
$ws.Cells.Item(31, 10).Value = 4230.3335
$ws.Cells.Item(31, 11).Value = 2436.2
$ws.Cells.Item(31, 12).Value = 4230.3335
$ws.Cells.Item(31, 13).Value = -2141.2
$ws.Cells.Item(31, 14).Value = -4820.3335

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3921
$ws.Cells.Item(34, 9).Value = 2436.2
$ws.Cells.Item(34, 10).Value = 4230.3335
$ws.Cells.Item(34, 11).Value = 2436.2
$ws.Cells.Item(34, 12).Value = 4230.3335
$ws.Cells.Item(34, 13).Value = -2234.2
$ws.Cells.Item(34, 14).Value = -4634.3335

# Sheet CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3045.5715
$ws.Cells.Item(62, 9).Value = 3090
$ws.Cells.Item(62, 10).Value = 2779
$ws.Cells.Item(62, 11).Value = 3090
$ws.Cells.Item(62, 12).Value = 2779
$ws.Cells.Item(62, 13).Value = -2466

# Sheet CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 3045.5715
$ws.Cells.Item(65, 9).Value = 3090
$ws.Cells.Item(65, 10).Value = 2779
$ws.Cells.Item(65, 11).Value = 15450
$ws.Cells.Item(65, 12).Value = 13895
$ws.Cells.Item(65, 13).Value = -12330

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3512.5483
$ws.Cells.Item(132, 9).Value = 3535.923
$ws.Cells.Item(132, 10).Value = 3391
$ws.Cells.Item(132, 11).Value = 10607.769
$ws.Cells.Item(132, 12).Value = 10173
$ws.Cells.Item(132, 13).Value = -8077.769
$ws.Cells.Item(132, 14).Value = -15233

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2233
$ws.Cells.Item(134, 9).Value = 2233
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 6699
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1697.9375
$ws.Cells.Item(5, 9).Value = 1610.5
$ws.Cells.Item(5, 10).Value = 1750.4
$ws.Cells.Item(5, 11).Value = 4831.5
$ws.Cells.Item(5, 12).Value = 5251.200000000001
$ws.Cells.Item(5, 13).Value = -4719.5
$ws.Cells.Item(5, 14).Value = -5475.200000000001

# Sheet CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 9974.5
$ws.Cells.Item(56, 9).Value = 9974.5
$ws.Cells.Item(56, 11).Value = 9974.5
$ws.Cells.Item(56, 13).Value = -9444.5

# Sheet CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 682
$ws.Cells.Item(68, 9).Value = 501
$ws.Cells.Item(68, 10).Value = 772.5
$ws.Cells.Item(68, 11).Value = 1503
$ws.Cells.Item(68, 12).Value = 2317.5
$ws.Cells.Item(68, 13).Value = -692
$ws.Cells.Item(68, 14).Value = -3939.5

# Sheet CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 682
$ws.Cells.Item(71, 9).Value = 501
$ws.Cells.Item(71, 10).Value = 772.5
$ws.Cells.Item(71, 11).Value = 4509
$ws.Cells.Item(71, 12).Value = 6952.5
$ws.Cells.Item(71, 13).Value = -453
$ws.Cells.Item(71, 14).Value = -15064.5

# Sheet CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 2207.4
$ws.Cells.Item(87, 9).Value = 2257.5
$ws.Cells.Item(87, 11).Value = 6772.5
$ws.Cells.Item(87, 13).Value = -5524.5

# Sheet CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 2207.4
$ws.Cells.Item(90, 9).Value = 2257.5
$ws.Cells.Item(90, 11).Value = 20317.5
$ws.Cells.Item(90, 13).Value = -14077.5

# Sheet CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1805.25
$ws.Cells.Item(113, 10).Value = 1937.5714
$ws.Cells.Item(113, 12).Value = 5812.7142
$ws.Cells.Item(113, 14).Value = -10152.7142

# Sheet CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 3141.3
$ws.Cells.Item(117, 10).Value = 3510.5
$ws.Cells.Item(117, 12).Value = 10531.5
$ws.Cells.Item(117, 14).Value = -17415.5

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 3333
$ws.Cells.Item(131, 9).Value = 3333
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 9999
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).ClearContents()

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 1697.9375
$ws.Cells.Item(135, 9).Value = 1610.5
$ws.Cells.Item(135, 10).Value = 1750.4
$ws.Cells.Item(135, 11).Value = 14494.5
$ws.Cells.Item(135, 12).Value = 15753.6
$ws.Cells.Item(135, 13).Value = -11959.5
$ws.Cells.Item(135, 14).Value = -20823.6

# Sheet GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 10046.25
$ws.Cells.Item(2, 9).Value = 26
$ws.Cells.Item(2, 11).Value = 26
$ws.Cells.Item(2, 13).Value = 87

# Sheet GSM row 96
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(96, 8).Value = 56666
$ws.Cells.Item(96, 10).Value = 56666
$ws.Cells.Item(96, 12).Value = 56666
$ws.Cells.Item(96, 14).Value = -62158

# Sheet GSM row 127
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(127, 8).Value = 92074.5
$ws.Cells.Item(127, 10).Value = 92074.5
$ws.Cells.Item(127, 12).Value = 92074.5
$ws.Cells.Item(127, 14).Value = -101994.5

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 7813.852
$ws.Cells.Item(46, 9).Value = 3149.5
$ws.Cells.Item(46, 10).Value = 8187
$ws.Cells.Item(46, 11).Value = 3149.5
$ws.Cells.Item(46, 12).Value = 8187
$ws.Cells.Item(46, 13).Value = -2961.5
$ws.Cells.Item(46, 14).Value = -8563

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4362.9067
$ws.Cells.Item(132, 9).Value = 4197.9033
$ws.Cells.Item(132, 10).Value = 4789.1665
$ws.Cells.Item(132, 11).Value = 12593.7099
$ws.Cells.Item(132, 12).Value = 14367.4995
$ws.Cells.Item(132, 13).Value = -10063.7099
$ws.Cells.Item(132, 14).Value = -19427.4995

# Sheet LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 12189.909
$ws.Cells.Item(136, 9).Value = 10012.571
$ws.Cells.Item(136, 10).Value = 16000.25
$ws.Cells.Item(136, 11).Value = 30037.713
$ws.Cells.Item(136, 12).Value = 48000.75
$ws.Cells.Item(136, 13).Value = -27487.713
$ws.Cells.Item(136, 14).Value = -53100.75

# Sheet WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 494.57144
$ws.Cells.Item(107, 9).Value = 352.5
$ws.Cells.Item(107, 10).Value = 684
$ws.Cells.Item(107, 11).Value = 1057.5
$ws.Cells.Item(107, 12).Value = 2052
$ws.Cells.Item(107, 13).Value = 862.5
$ws.Cells.Item(107, 14).Value = -5892

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1674.8125
$ws.Cells.Item(132, 9).Value = 1586.5
$ws.Cells.Item(132, 10).Value = 2999.5
$ws.Cells.Item(132, 11).Value = 4759.5
$ws.Cells.Item(132, 12).Value = 8998.5
$ws.Cells.Item(132, 13).Value = -2229.5
$ws.Cells.Item(132, 14).Value = -14058.5
